$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-9 from 45233 to 45243
$ws.Range("C2:C9").Value = 45243
